# Applies cryptos list update per commit "Updated cryptos list on Sat Feb 24 21:44:16 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.484.39"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "2.981.47"
$ws.Range("E3").Value = "  +1.21%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "381.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.84%  "
$ws.Range("E7").Value = "  +0.79%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.591"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.61"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.137"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("D13").Value = "3.458.97"
$ws.Range("E13").Value = "  +1.68%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.83%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "18.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("D16").Value = "2.995.45"
$ws.Range("E16").Value = "  +1.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "11.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.994"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("D19").Value = "51.512.27"
$ws.Range("E19").Value = "  +1.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.75%  "
$ws.Range("D22").Value = "0.0₃0962"
$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("E25").Value = "  +2.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.18%  "
$ws.Range("E28").Value = "  +3.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.41"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "51.38"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.85%  "
$ws.Range("E35").Value = "  +0.63%  "
$ws.Range("E36").Value = "  +0.96%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.28"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.94"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.67%  "
$ws.Range("E40").Value = "  +4.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.116"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +11.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "126.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.271"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.41%  "
$ws.Range("D49").Value = "2.021.63"
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("D50").Value = "3.278.36"
$ws.Range("E50").Value = "  +1.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0328"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.22%  "
